# Cebollín (Terminal Hortofrutícola Agro Chillán) weekly refresh.
# A new weekly price observation is inserted at row 202, pushing the
# previously-recorded rows 202-217 down to 203-218 (dimension grows from
# A1:R217 to A1:R218).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 202; Excel shifts rows 202:217 down to 203:218
# and carries the row-below's formatting (including the date style on
# column D) up into the freshly inserted row, same as native Excel.
$ws.Rows("202:202").Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A202").Value = 7
$ws.Range("B202").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C202").Value = "Ñuble"
$ws.Range("D202").Value = 45212
$ws.Range("E202").Value = 16
$ws.Range("F202").Value = 100112037
$ws.Range("G202").Value = "Cebollín"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 150
$ws.Range("K202").Value = 5000
$ws.Range("L202").Value = 5000
$ws.Range("M202").Value = 5000
$ws.Range("N202").Value = "$/paquete 36 unidades"
$ws.Range("O202").Value = "Provincia de Diguillín"
$ws.Range("P202").Value = 139
$ws.Range("Q202").Value = 36
$ws.Range("R202").Value = "Hortaliza"
